# Applies the commit "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
# Summary of the real (semantic) content changes found in the target diff:
#   1. E11 (VALOR MORA)      : 10800000 -> 10900000
#   2. F13 (Cant. Periodos)  : 108 -> 109
#   3. The "Periodo Mora" table (col E, rows 16..) lists every period from
#      1607 (jul-2016) through 2506 (jun-2025) ascending. The new workbook
#      drops nothing from that span but ADDS period 2507 (jul-2025) and
#      re-lists everything in DESCENDING order (newest first), so the table
#      grows by one row (16..123 -> 16..124).
#   4. Row 123 (previously the last/bottom-bordered row of the table) must
#      become a normal interior row, and a new row 124 becomes the
#      bottom-bordered closing row of the table.
#   5. The signature block (two label rows just below the table) shifts
#      down by one row: old row 128/129 -> new row 129/130. Text content is
#      unchanged, only its row position moves.
# All the other <v>NN</v> index churn visible in the raw xml diff is just
# sharedStrings renumbering caused by inserting the single new "2507"
# string - the actual label text on every one of those cells is identical
# before/after, so none of those cells need to be touched here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Simple value edits
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 10900000
$ws.Range("F13").Value = 109

# ---------------------------------------------------------------------
# 2) Make room for the extra period: insert a new row at 124 (this also
#    carries the old rows 124.. (the blank rows + the signature block
#    128/129) one row further down, landing on 129/130).
# ---------------------------------------------------------------------
$ws.Rows.Item(124).Insert()

# Give the brand-new row 124 the table's "closing row" look (bottom
# border) by copying the format that row 123 still has at this point.
$ws.Range("B123:J123").Copy()
$ws.Range("B124:J124").PasteSpecial(-4122)   # xlPasteFormats

# Row 123 now becomes an ordinary interior row of the table, so give it
# the same look as every other data row (e.g. row 16).
$ws.Range("B16:J16").Copy()
$ws.Range("B123:J123").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Re-populate column E (Periodo Mora) for every row of the table,
#    rows 16..124, newest period first (descending), including the new
#    period 2507.
# ---------------------------------------------------------------------
$periods = @()
for ($y = 25; $y -ge 16; $y--) {
    for ($m = 12; $m -ge 1; $m--) {
        if ($y -eq 16 -and $m -lt 7) { continue }
        if ($y -eq 25 -and $m -gt 7) { continue }
        $periods += ('{0:D2}{1:D2}' -f $y, $m)
    }
}

for ($i = 0; $i -lt $periods.Count; $i++) {
    $ws.Cells.Item(16 + $i, 5).Value = $periods[$i]
}

# The new row 124 needs the same "template" values every other table row
# carries in columns B/C/D/F/G (only E - the period - varies row to row).
$ws.Range("B124").Value = $ws.Range("B123").Value
$ws.Range("C124").Value = $ws.Range("C123").Value
$ws.Range("D124").Value = $ws.Range("D123").Value
$ws.Range("F124").Value = $ws.Range("F123").Value
$ws.Range("G124").Value = $ws.Range("G123").Value

Write-Output "Populated $($periods.Count) periods from $($periods[0]) down to $($periods[$periods.Count-1])"
